# Swap the presentation's theme colour scheme back to the stock
# "Office" palette (it had been re-coloured to the "Integral" theme
# palette). This is done through the Theme Colors object model
# (Slide.ThemeColorScheme / SlideRange.ThemeColorScheme), which is the
# supported COM surface for editing the colours stored in the design's
# theme part (ppt/theme/theme2.xml) - the same part backing every slide,
# slide layout and slide master in this deck.
#
# ThemeColorScheme items are ordered exactly like <a:clrScheme> in the
# OOXML theme part:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
#
# RGB is supplied the usual VBA/COM way: R + G*256 + B*65536.

$p = $ppt.ActivePresentation

function Set-ThemeRgb($tcs, $index, $r, $g, $b) {
    $tcs.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Use the first slide as the entry point into the shared theme colour
# scheme - any slide/slide range works since they all resolve to the
# same underlying theme part.
$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRgb $tcs 1  0x00 0x00 0x00   # dk1      -> 000000
Set-ThemeRgb $tcs 2  0xFF 0xFF 0xFF   # lt1      -> FFFFFF
Set-ThemeRgb $tcs 3  0x44 0x54 0x6A   # dk2      -> 44546A
Set-ThemeRgb $tcs 4  0xE7 0xE6 0xE6   # lt2      -> E7E6E6
Set-ThemeRgb $tcs 5  0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
Set-ThemeRgb $tcs 6  0xED 0x7D 0x31   # accent2  -> ED7D31
Set-ThemeRgb $tcs 7  0xA5 0xA5 0xA5   # accent3  -> A5A5A5
Set-ThemeRgb $tcs 8  0xFF 0xC0 0x00   # accent4  -> FFC000
Set-ThemeRgb $tcs 9  0x44 0x72 0xC4   # accent5  -> 4472C4
Set-ThemeRgb $tcs 10 0x70 0xAD 0x47   # accent6  -> 70AD47
Set-ThemeRgb $tcs 11 0x05 0x63 0xC1   # hlink    -> 0563C1
Set-ThemeRgb $tcs 12 0x95 0x4F 0x72   # folHlink -> 954F72
